$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 13752.667
$ws.Range("I18").Value = 5666.6665
$ws.Range("J18").Value = 17795.666
$ws.Range("K18").Value = 5666.6665
$ws.Range("L18").Value = 17795.666
$ws.Range("M18").Value = -5382.6665
$ws.Range("N18").Value = -18363.666
$ws.Range("H38").Value = 579.8182
$ws.Range("I38").Value = 437.8
$ws.Range("K38").Value = 1313.4
$ws.Range("M38").Value = -941.4000000000001
$ws.Range("H43").Value = 849.4545000000001
$ws.Range("J43").Value = 899.44446
$ws.Range("L43").Value = 899.44446
$ws.Range("N43").Value = -1037.44446
$ws.Range("H69").Value = 7450
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7450
$ws.Range("K69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("M69").Value = 22350
$ws.Range("N69").Value = -24098
$ws.Range("H70").Value = 10623.637
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 11586
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 34758
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -35298
$ws.Range("H72").Value = 7450
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7450
$ws.Range("K72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("M72").Value = 67050
$ws.Range("N72").Value = -75786
$ws.Range("H73").Value = 10623.637
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 11586
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 34758
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -36630
$ws.Range("H131").Value = 1560.125
$ws.Range("I131").Value = 834
$ws.Range("J131").Value = 3157.6
$ws.Range("K131").Value = 2502
$ws.Range("L131").Value = 9472.799999999999
$ws.Range("M131").Value = 2538
$ws.Range("N131").Value = -19552.8
$ws.Range("H132").Value = 559.89655
$ws.Range("I132").Value = 508.68674
$ws.Range("K132").Value = 1526.06022
$ws.Range("M132").Value = 1003.93978
$ws.Range("H138").Value = 2478.8774
$ws.Range("I138").Value = 2585.9614
$ws.Range("J138").Value = 2357.8262
$ws.Range("K138").Value = 7757.8842
$ws.Range("L138").Value = 7073.4786
$ws.Range("M138").Value = -2617.8842
$ws.Range("N138").Value = -17353.4786
$ws.Range("H141").Value = 4043.077
$ws.Range("I141").Value = 1462.5
$ws.Range("K141").Value = 4387.5
$ws.Range("M141").Value = 792.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4193
$ws.Range("I32").Value = 3120.8572
$ws.Range("K32").Value = 3120.8572
$ws.Range("M32").Value = -2833.8572
$ws.Range("H39").Value = 4999
$ws.Range("I39").Value = 4999
$ws.Range("K39").Value = 4999
$ws.Range("M39").Value = -4479
$ws.Range("H95").Value = 65103.5
$ws.Range("J95").Value = 65103.5
$ws.Range("L95").Value = 65103.5
$ws.Range("N95").Value = -70595.5
$ws.Range("H97").Value = 1175.75
$ws.Range("I97").Value = 1175.75
$ws.Range("K97").Value = 1175.75
$ws.Range("M97").Value = -679.75
$ws.Range("H101").Value = 42483
$ws.Range("J101").Value = 42483
$ws.Range("L101").Value = 42483
$ws.Range("N101").Value = -48973
$ws.Range("H114").Value = 9249
$ws.Range("J114").Value = 9249
$ws.Range("L114").Value = 9249
$ws.Range("N114").Value = -17927
$ws.Range("H132").Value = 2573.5715
$ws.Range("I132").Value = 1848.1111
$ws.Range("K132").Value = 5544.3333
$ws.Range("M132").Value = -3014.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2500
$ws.Range("I7").Value = 2500
$ws.Range("K7").Value = 2500
$ws.Range("M7").Value = -2387
$ws.Range("H20").Value = 1645.3
$ws.Range("I20").Value = 1506
$ws.Range("K20").Value = 1506
$ws.Range("M20").Value = -1259
$ws.Range("H81").Value = 19396.5
$ws.Range("J81").Value = 19396.5
$ws.Range("L81").Value = 19396.5
$ws.Range("N81").Value = -21518.5
$ws.Range("H84").Value = 19396.5
$ws.Range("J84").Value = 19396.5
$ws.Range("L84").Value = 58189.5
$ws.Range("N84").Value = -68797.5
$ws.Range("H100").Value = 25999.666
$ws.Range("J100").Value = 25999.666
$ws.Range("L100").Value = 25999.666
$ws.Range("N100").Value = -28163.666
$ws.Range("H130").Value = 36665.332
$ws.Range("J130").Value = 36665.332
$ws.Range("L130").Value = 36665.332
$ws.Range("N130").Value = -46705.332

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3031.2856
$ws.Range("I31").Value = 3702.4
$ws.Range("J31").Value = 2658.4443
$ws.Range("K31").Value = 3702.4
$ws.Range("L31").Value = 2658.4443
$ws.Range("M31").Value = -3407.4
$ws.Range("N31").Value = -3248.4443
$ws.Range("H34").Value = 3031.2856
$ws.Range("I34").Value = 3702.4
$ws.Range("J34").Value = 2658.4443
$ws.Range("K34").Value = 3702.4
$ws.Range("L34").Value = 2658.4443
$ws.Range("M34").Value = -3500.4
$ws.Range("N34").Value = -3062.4443
$ws.Range("H43").Value = 23333.334
$ws.Range("J43").Value = 23333.334
$ws.Range("L43").Value = 23333.334
$ws.Range("N43").Value = -23701.334
$ws.Range("H95").Value = 26236
$ws.Range("J95").Value = 26236
$ws.Range("L95").Value = 26236
$ws.Range("N95").Value = -31728
$ws.Range("H101").Value = 23333.334
$ws.Range("J101").Value = 23333.334
$ws.Range("L101").Value = 23333.334
$ws.Range("N101").Value = -29823.334
$ws.Range("H132").Value = 2001.7297
$ws.Range("J132").Value = 3483.5715
$ws.Range("L132").Value = 10450.7145
$ws.Range("N132").Value = -15510.7145
$ws.Range("H134").Value = 1835.6875
$ws.Range("I134").Value = 1577.3448
$ws.Range("K134").Value = 4732.0344
$ws.Range("M134").Value = -2197.0344

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 6500
$ws.Range("J70").Value = 6500
$ws.Range("L70").Value = 19500
$ws.Range("N70").Value = -20130
$ws.Range("H73").Value = 6500
$ws.Range("J73").Value = 6500
$ws.Range("L73").Value = 19500
$ws.Range("N73").Value = -21684
$ws.Range("H122").Value = 772
$ws.Range("I122").Value = 381.25
$ws.Range("J122").Value = 1032.5
$ws.Range("K122").Value = 3431.25
$ws.Range("L122").Value = 9292.5
$ws.Range("M122").Value = -981.25
$ws.Range("N122").Value = -14192.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1866.6666
$ws.Range("I80").Value = 1866.6666
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1866.6666
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = -868.6666
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 1866.6666
$ws.Range("I83").Value = 1866.6666
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 9333.333000000001
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = -4341.333000000001
$ws.Range("N83").ClearContents()
$ws.Range("H92").Value = 22537
$ws.Range("J92").Value = 22537
$ws.Range("L92").Value = 22537
$ws.Range("N92").Value = -26281
$ws.Range("H98").Value = 29401
$ws.Range("J98").Value = 29401
$ws.Range("L98").Value = 29401
$ws.Range("N98").Value = -35391
$ws.Range("H102").Value = 1481.0975
$ws.Range("J102").Value = 1239.4
$ws.Range("L102").Value = 1239.4
$ws.Range("N102").Value = -4483.4
$ws.Range("H132").Value = 2568046.5
$ws.Range("I132").Value = 6413087
$ws.Range("K132").Value = 19239261
$ws.Range("M132").Value = -19236731

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 1000000000
$ws.Range("I47").Value = 1000000000
$ws.Range("K47").Value = 1000000000
$ws.Range("M47").Value = -999999510
$ws.Range("H52").Value = 1000000000
$ws.Range("I52").Value = 1000000000
$ws.Range("K52").Value = 1000000000
$ws.Range("M52").Value = -999999767
$ws.Range("H94").Value = 36388.168
$ws.Range("J94").Value = 36388.168
$ws.Range("L94").Value = 36388.168
$ws.Range("N94").Value = -37740.168
$ws.Range("H122").Value = 1743.72
$ws.Range("I122").Value = 1009.4545
$ws.Range("J122").Value = 2320.6428
$ws.Range("K122").Value = 3028.3635
$ws.Range("L122").Value = 6961.928400000001
$ws.Range("M122").Value = -578.3635000000004
$ws.Range("N122").Value = -11861.9284
$ws.Range("H132").Value = 2041.5454
$ws.Range("I132").Value = 1343.1666
$ws.Range("K132").Value = 4029.4998
$ws.Range("M132").Value = -1499.4998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 3000
$ws.Range("J61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3584
$ws.Range("H63").Value = 25999.666
$ws.Range("J63").Value = 25999.666
$ws.Range("L63").Value = 25999.666
$ws.Range("N63").Value = -27247.666
$ws.Range("H66").Value = 25999.666
$ws.Range("J66").Value = 25999.666
$ws.Range("L66").Value = 77998.99800000001
$ws.Range("N66").Value = -84238.99800000001
$ws.Range("H69").Value = 13824.25
$ws.Range("J69").Value = 13824.25
$ws.Range("L69").Value = 13824.25
$ws.Range("N69").Value = -15322.25
$ws.Range("H72").Value = 13824.25
$ws.Range("J72").Value = 13824.25
$ws.Range("L72").Value = 41472.75
$ws.Range("N72").Value = -48960.75
$ws.Range("H105").Value = 43293.668
$ws.Range("J105").Value = 43293.668
$ws.Range("L105").Value = 43293.668
$ws.Range("N105").Value = -50281.668
$ws.Range("H122").Value = 158735.4
$ws.Range("I122").Value = 158735.4
$ws.Range("K122").Value = 476206.2
$ws.Range("M122").Value = -473756.2
$ws.Range("H132").Value = 2849.6667
$ws.Range("I132").Value = 1775
$ws.Range("K132").Value = 5325
$ws.Range("M132").Value = -2795
